$d = $word.ActiveDocument

# 1. Replace the placeholder title text with the real lab title.
$rng = $d.Content
$null = $rng.Find.Execute("Назва", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Моделювання систем масового обслуговування з одним обслуговуючим приладом та чергою", 2)

# 2. Update the H1 / H1 Char style pair: stop auto-redefining from usage and
#    render the heading in all caps.
$h1 = $d.Styles("H1")
$h1.AutomaticallyUpdate = $false
$h1.Font.AllCaps = $true

$h1Char = $d.Styles("H1Char")
$h1Char.Font.AllCaps = $true
